# Update "want-to-go" counts (column F) across the four worksheets to match
# the newly generated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 387
$ws1.Range("F4").Value  = 155
$ws1.Range("F5").Value  = 1318
$ws1.Range("F6").Value  = 229
$ws1.Range("F7").Value  = 2509
$ws1.Range("F8").Value  = 914
$ws1.Range("F9").Value  = 18691
$ws1.Range("F11").Value = 1931
$ws1.Range("F12").Value = 669
$ws1.Range("F14").Value = 334
$ws1.Range("F15").Value = 606
$ws1.Range("F17").Value = 204
$ws1.Range("F20").Value = 23
$ws1.Range("F23").Value = 105

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value  = 111
$ws2.Range("F10").Value = 226

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5896
$ws3.Range("F3").Value = 573
$ws3.Range("F4").Value = 557

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 5896
$ws4.Range("F4").Value  = 573
$ws4.Range("F5").Value  = 557
$ws4.Range("F6").Value  = 387
$ws4.Range("F8").Value  = 155
$ws4.Range("F10").Value = 1318
$ws4.Range("F12").Value = 229
$ws4.Range("F15").Value = 2509
$ws4.Range("F16").Value = 914
$ws4.Range("F17").Value = 18691
$ws4.Range("F21").Value = 111
$ws4.Range("F22").Value = 226
$ws4.Range("F23").Value = 226
$ws4.Range("F24").Value = 1931
$ws4.Range("F25").Value = 669
$ws4.Range("F27").Value = 334
$ws4.Range("F28").Value = 606
$ws4.Range("F30").Value = 204
$ws4.Range("F36").Value = 23
$ws4.Range("F41").Value = 105

$wb.Save()
